$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K99:L99").ClearFormats()
